$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 (Item #3, Designator J3) part substitution: Omron -> CONNFLY Elec connector
$ws.Range("D9").Value = "CONNFLY Elec"
$ws.Range("E9").Value = "DS1037-15FNAKT74-0CC"
$ws.Range("I9").Value = "LCSC Part: C77835"
